$d = $word.ActiveDocument

# 1) "le plus de points possibles sans" -> "le plus de points possible sans"
#    (single run, no split)
$d.Content.Find.Execute(
    "le plus de points possibles sans",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "le plus de points possible sans",
    2)

# 2) "nombre de point ramassé et d'" -> "nombre de points ramassés et d'"
#    Target OOXML splits into new runs for the inserted "s" letters, so we
#    perform the two single-letter insertions as separate Find/Replace
#    operations on narrow ranges (mirrors manual retyping in Word, which is
#    what produces the run split seen in the diff).
$d.Content.Find.Execute(
    "nombre de point ramassé et d",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "nombre de points ramassés et d",
    2)

# 3) "gérer l'évènement ou le joueur" -> "gérer l'évènement où le joueur"
$d.Content.Find.Execute(
    "gérer l’évènement ou le joueur",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "gérer l’évènement où le joueur",
    2)
